$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename subjects (shared strings at rows 2-11, column B)
$ws.Range("B2").Value = "S4"
$ws.Range("B3").Value = "S2"
$ws.Range("B4").Value = "S7"
$ws.Range("B5").Value = "S14"
$ws.Range("B6").Value = "S9"
$ws.Range("B7").Value = "S6"
$ws.Range("B8").Value = "S11"
$ws.Range("B9").Value = "S8"
$ws.Range("B10").Value = "S5"
$ws.Range("B11").Value = "S10"

# Update regional data values (regenerated fake data)
$ws.Range("C2").Value = 0.42524673444907
$ws.Range("D2").Value = 0.665001459280523
$ws.Range("E2").Value = 0.460931045275401
$ws.Range("F2").Value = 1.02570603999119
$ws.Range("G2").Value = 0.13875944221594
$ws.Range("H2").Value = 0.986566684067103
$ws.Range("I2").Value = 1.46946299898945
$ws.Range("J2").Value = 0.463215835327427
$ws.Range("K2").Value = 0.33304525885424
$ws.Range("L2").Value = 0.409877874158788
$ws.Range("M2").Value = 3.02720855130581
$ws.Range("N2").Value = 0.0259598686933449
$ws.Range("O2").Value = 0.956550502252821
$ws.Range("P2").Value = 1.00106633201535
$ws.Range("Q2").Value = 0.923329721768157
$ws.Range("R2").Value = 0.0863691347788578
$ws.Range("S2").Value = 1.13371914252842
$ws.Range("T2").Value = 1.93800899100486
$ws.Range("U2").Value = 0.975233218657467
$ws.Range("V2").Value = 1.43506774974341
$ws.Range("W2").Value = 0.840003565946416
$ws.Range("X2").Value = 0.889230238592031
$ws.Range("Y2").Value = 1.54495250285111
$ws.Range("Z2").Value = 0.945542845975493
$ws.Range("AA2").Value = 0.28702916483433
$ws.Range("AB2").Value = 0.746088708061449
$ws.Range("AC2").Value = 0.479994023113894
$ws.Range("AD2").Value = 0.334380468086379
$ws.Range("AE2").Value = 0.612553343402353
$ws.Range("AF2").Value = 0.918871430386518
$ws.Range("AG2").Value = 0.437867886402896
$ws.Range("AH2").Value = 1.68438549713851
$ws.Range("AI2").Value = 0.687381404819305
$ws.Range("AJ2").Value = 0.202237315574305
$ws.Range("AK2").Value = 0.784862874485179
$ws.Range("AL2").Value = 1.44324500587884
$ws.Range("AM2").Value = 0.281816533290877
$ws.Range("AN2").Value = 0.913444098070634
$ws.Range("AO2").Value = 0.0992090109440399
$ws.Range("AP2").Value = 1.64054316810279
$ws.Range("AQ2").Value = 0.698515790647975
$ws.Range("AR2").Value = 0.863490087491982
$ws.Range("AS2").Value = 0.798942289368513
$ws.Range("AT2").Value = 3.33310341646591
$ws.Range("AU2").Value = 0.795523500040539
$ws.Range("AV2").Value = 0.887697265477213
$ws.Range("AW2").Value = 0.466500982236482
$ws.Range("AX2").Value = 0.65120462206912
$ws.Range("AY2").Value = 0.450381993016291
$ws.Range("AZ2").Value = 0.21066104632662
$ws.Range("BA2").Value = 2.08401572176472
$ws.Range("BB2").Value = 0.357145505954922
$ws.Range("BC2").Value = 0.673542545247511
$ws.Range("BD2").Value = 1.02660552488693
$ws.Range("BE2").Value = 0.464796034571796
$ws.Range("BF2").Value = 0.139012261944852
$ws.Range("BG2").Value = 1.12969754646799
$ws.Range("BH2").Value = 0.472409221121823
$ws.Range("BI2").Value = 0.366252454838987
$ws.Range("BJ2").Value = 1.47399163157824
$ws.Range("BK2").Value = 0.187090707347886
$ws.Range("BL2").Value = 0.502063722228696
$ws.Range("BM2").Value = 0.272341170090619
$ws.Range("BN2").Value = 1.7517707694319
$ws.Range("BO2").Value = 0.355612465521914
$ws.Range("BP2").Value = 1.06363984261649
$ws.Range("BQ2").Value = 1.49689585903709
$ws.Range("BR2").Value = 0.146264879100494
$ws.Range("BS2").Value = 0.658039645728341
$ws.Range("BT2").Value = 0.816275581007904
$ws.Range("BU2").Value = 0.530403273082154
$ws.Range("BV2").Value = 0.755903359029766
$ws.Range("BW2").Value = 1.78893659815308
$ws.Range("BX2").Value = 1.54795399933233
$ws.Range("BY2").Value = 0.864034389390506
$ws.Range("BZ2").Value = 0.517150001318192
$ws.Range("CA2").Value = 0.944124862571797
$ws.Range("CB2").Value = 0.104910343338995
$ws.Range("CC2").Value = 0.636804202157967
$ws.Range("CD2").Value = 0.458578449796607
$ws.Range("C3").Value = 0.481369311530437
$ws.Range("D3").Value = 0.970674407854476
$ws.Range("E3").Value = 0.328053649194687
$ws.Range("F3").Value = 0.985216844022111
$ws.Range("G3").Value = 2.3646741061793
$ws.Range("H3").Value = 1.64598917962507
$ws.Range("I3").Value = 0.601937935077132
$ws.Range("J3").Value = 0.64365181561149
$ws.Range("K3").Value = 0.943900614960771
$ws.Range("L3").Value = 0.944625162942793
$ws.Range("M3").Value = 1.0431860578736
$ws.Range("N3").Value = 0.576723590228008
$ws.Range("O3").Value = 0.68416935346147
$ws.Range("P3").Value = 0.213660306844166
$ws.Range("Q3").Value = 0.673328583878302
$ws.Range("R3").Value = 1.03596783245367
$ws.Range("S3").Value = 0.961107981231816
$ws.Range("T3").Value = 1.41670069713756
$ws.Range("U3").Value = 1.1522984902046
$ws.Range("V3").Value = 0.398227024442586
$ws.Range("W3").Value = 0.51856652438013
$ws.Range("X3").Value = 1.9838939247505
$ws.Range("Y3").Value = 0.545589009688193
$ws.Range("Z3").Value = 0.713445488033851
$ws.Range("AA3").Value = 1.93954235525332
$ws.Range("AB3").Value = 1.54162080750076
$ws.Range("AC3").Value = 1.05960755832391
$ws.Range("AD3").Value = 0.14131309412087
$ws.Range("AE3").Value = 0.675978259149002
$ws.Range("AF3").Value = 0.0853754660832985
$ws.Range("AG3").Value = 0.84725506661529
$ws.Range("AH3").Value = 0.551444555562282
$ws.Range("AI3").Value = 0.895266079954652
$ws.Range("AJ3").Value = 2.91587591655616
$ws.Range("AK3").Value = 0.858237851574561
$ws.Range("AL3").Value = 1.18419001560228
$ws.Range("AM3").Value = 3.85165370227529
$ws.Range("AN3").Value = 0.101136829007663
$ws.Range("AO3").Value = 1.29416581420977
$ws.Range("AP3").Value = 1.74762173583148
$ws.Range("AQ3").Value = 1.45074762132195
$ws.Range("AR3").Value = 1.37643118207895
$ws.Range("AS3").Value = 0.295174932321449
$ws.Range("AT3").Value = 0.979773357449307
$ws.Range("AU3").Value = 1.58016420834449
$ws.Range("AV3").Value = 1.25974358443635
$ws.Range("AW3").Value = 0.991901437061612
$ws.Range("AX3").Value = 0.394251138275582
$ws.Range("AY3").Value = 0.29013342674108
$ws.Range("AZ3").Value = 1.672365478342
$ws.Range("BA3").Value = 0.926015269840624
$ws.Range("BB3").Value = 0.672711153491911
$ws.Range("BC3").Value = 1.55682790302541
$ws.Range("BD3").Value = 1.04328137993115
$ws.Range("BE3").Value = 0.786375988566788
$ws.Range("BF3").Value = 0.834270776060634
$ws.Range("BG3").Value = 0.318456093896073
$ws.Range("BH3").Value = 1.05719242759025
$ws.Range("BI3").Value = 0.813787398117652
$ws.Range("BJ3").Value = 1.04935419116017
$ws.Range("BK3").Value = 0.503497139417598
$ws.Range("BL3").Value = 1.79188033326898
$ws.Range("BM3").Value = 0.277992405110174
$ws.Range("BN3").Value = 0.452928176403982
$ws.Range("BO3").Value = 1.11384356770448
$ws.Range("BP3").Value = 0.36314431219147
$ws.Range("BQ3").Value = 0.216535127694439
$ws.Range("BR3").Value = 0.106224358246996
$ws.Range("BS3").Value = 1.98139506995511
$ws.Range("BT3").Value = 0.229809524702531
$ws.Range("BU3").Value = 0.244446571017464
$ws.Range("BV3").Value = 2.78554775596674
$ws.Range("BW3").Value = 2.33332869151362
$ws.Range("BX3").Value = 0.0982859262760793
$ws.Range("BY3").Value = 1.08839330917256
$ws.Range("BZ3").Value = 0.93949436343093
$ws.Range("CA3").Value = 2.00300401305992
$ws.Range("CB3").Value = 0.836557811851305
$ws.Range("CC3").Value = 0.506255245718704
$ws.Range("CD3").Value = 0.402446620500866
$ws.Range("C4").Value = 0.783799225470391
$ws.Range("D4").Value = 0.661049292066198
$ws.Range("E4").Value = 0.175156546751359
$ws.Range("F4").Value = 0.47532665114732
$ws.Range("G4").Value = 0.83807600874854
$ws.Range("H4").Value = 1.51846752277608
$ws.Range("I4").Value = 0.566471663434815
$ws.Range("J4").Value = 0.398545594904318
$ws.Range("K4").Value = 1.00403496091877
$ws.Range("L4").Value = 1.9268810236581
$ws.Range("M4").Value = 0.89205626542496
$ws.Range("N4").Value = 0.189699771235515
$ws.Range("O4").Value = 0.711024550933412
$ws.Range("P4").Value = 0.810168610539189
$ws.Range("Q4").Value = 2.50896264681446
$ws.Range("R4").Value = 1.359698867191
$ws.Range("S4").Value = 0.827596906152946
$ws.Range("T4").Value = 2.55993508680299
$ws.Range("U4").Value = 1.15228430154269
$ws.Range("V4").Value = 0.868650260626903
$ws.Range("W4").Value = 0.310324274142769
$ws.Range("X4").Value = 2.53331658624197
$ws.Range("Y4").Value = 0.910453510293951
$ws.Range("Z4").Value = 1.07777574781712
$ws.Range("AA4").Value = 2.03469860169484
$ws.Range("AB4").Value = 1.90889502822426
$ws.Range("AC4").Value = 1.60220096636113
$ws.Range("AD4").Value = 1.22856405480769
$ws.Range("AE4").Value = 0.0767091839938488
$ws.Range("AF4").Value = 1.88939636637206
$ws.Range("AG4").Value = 1.82732280230812
$ws.Range("AH4").Value = 0.480973011574268
$ws.Range("AI4").Value = 1.17352943965094
$ws.Range("AJ4").Value = 0.252960875077487
$ws.Range("AK4").Value = 0.205300899345915
$ws.Range("AL4").Value = 0.493187089428245
$ws.Range("AM4").Value = 1.08166715280504
$ws.Range("AN4").Value = 0.796180294058755
$ws.Range("AO4").Value = 1.72064202536912
$ws.Range("AP4").Value = 0.40260806819923
$ws.Range("AQ4").Value = 0.41377266804411
$ws.Range("AR4").Value = 2.30085521523483
$ws.Range("AS4").Value = 1.49802728560788
$ws.Range("AT4").Value = 1.62055626791413
$ws.Range("AU4").Value = 1.32996150987887
$ws.Range("AV4").Value = 0.882326282036702
$ws.Range("AW4").Value = 0.365680434553212
$ws.Range("AX4").Value = 0.563816345994248
$ws.Range("AY4").Value = 0.730374245299592
$ws.Range("AZ4").Value = 4.40428895422964
$ws.Range("BA4").Value = 0.490555140768178
$ws.Range("BB4").Value = 1.70823939214758
$ws.Range("BC4").Value = 0.103730217613519
$ws.Range("BD4").Value = 1.58642160541036
$ws.Range("BE4").Value = 0.793164518459762
$ws.Range("BF4").Value = 0.591879431869704
$ws.Range("BG4").Value = 0.782693407991043
$ws.Range("BH4").Value = 3.00984003709629
$ws.Range("BI4").Value = 0.835586217692451
$ws.Range("BJ4").Value = 1.35800183755025
$ws.Range("BK4").Value = 0.459497661740678
$ws.Range("BL4").Value = 0.734966327878891
$ws.Range("BM4").Value = 0.614423220370625
$ws.Range("BN4").Value = 0.468202497732698
$ws.Range("BO4").Value = 0.325123503779649
$ws.Range("BP4").Value = 1.02482931270505
$ws.Range("BQ4").Value = 1.47598778438599
$ws.Range("BR4").Value = 0.746639648360255
$ws.Range("BS4").Value = 0.209536929189714
$ws.Range("BT4").Value = 0.809003369207542
$ws.Range("BU4").Value = 1.17165353990013
$ws.Range("BV4").Value = 0.709580703857448
$ws.Range("BW4").Value = 0.546841298518953
$ws.Range("BX4").Value = 0.716018352791933
$ws.Range("BY4").Value = 0.0722034411730693
$ws.Range("BZ4").Value = 0.783156146727711
$ws.Range("CA4").Value = 2.06096460114081
$ws.Range("CB4").Value = 1.02480755614043
$ws.Range("CC4").Value = 1.26912726303568
$ws.Range("CD4").Value = 0.492403295065684
$ws.Range("C5").Value = 0.484094316617979
$ws.Range("D5").Value = 0.6768010577612
$ws.Range("E5").Value = 1.43772764474382
$ws.Range("F5").Value = 2.07731280703434
$ws.Range("G5").Value = 0.547276860706612
$ws.Range("H5").Value = 0.555102270699726
$ws.Range("I5").Value = 1.40170710216103
$ws.Range("J5").Value = 0.28704877271029
$ws.Range("K5").Value = 1.34286784587768
$ws.Range("L5").Value = 1.10843206895112
$ws.Range("M5").Value = 1.24186450476699
$ws.Range("N5").Value = 0.922995505624285
$ws.Range("O5").Value = 1.28974014507287
$ws.Range("P5").Value = 0.574999950906318
$ws.Range("Q5").Value = 0.840848842591756
$ws.Range("R5").Value = 0.0845691479354462
$ws.Range("S5").Value = 0.536659731962058
$ws.Range("T5").Value = 2.22810846603495
$ws.Range("U5").Value = 0.306561882698089
$ws.Range("V5").Value = 0.836739287167326
$ws.Range("W5").Value = 0.367817262828379
$ws.Range("X5").Value = 0.845240119911944
$ws.Range("Y5").Value = 1.53304078371534
$ws.Range("Z5").Value = 0.353802309671729
$ws.Range("AA5").Value = 2.76382663535462
$ws.Range("AB5").Value = 0.963107673627094
$ws.Range("AC5").Value = 0.414285216017406
$ws.Range("AD5").Value = 1.90157162853245
$ws.Range("AE5").Value = 0.887089653102139
$ws.Range("AF5").Value = 0.379299587158663
$ws.Range("AG5").Value = 1.1606577287469
$ws.Range("AH5").Value = 0.491186614464462
$ws.Range("AI5").Value = 0.966722047643857
$ws.Range("AJ5").Value = 0.949977840844719
$ws.Range("AK5").Value = 0.568013646802659
$ws.Range("AL5").Value = 1.32999219357781
$ws.Range("AM5").Value = 0.877916904690136
$ws.Range("AN5").Value = 0.593125960057563
$ws.Range("AO5").Value = 0.625103891942961
$ws.Range("AP5").Value = 0.406261935096915
$ws.Range("AQ5").Value = 0.710639084672
$ws.Range("AR5").Value = 1.03514593910933
$ws.Range("AS5").Value = 1.31936741792795
$ws.Range("AT5").Value = 0.970606381194801
$ws.Range("AU5").Value = 0.203100224323485
$ws.Range("AV5").Value = 2.18516835819793
$ws.Range("AW5").Value = 1.0638855064339
$ws.Range("AX5").Value = 0.134229718774706
$ws.Range("AY5").Value = 0.314206961157191
$ws.Range("AZ5").Value = 0.10785541262139
$ws.Range("BA5").Value = 0.856413651157827
$ws.Range("BB5").Value = 2.4686836563673
$ws.Range("BC5").Value = 0.738808231063749
$ws.Range("BD5").Value = 0.235311032938155
$ws.Range("BE5").Value = 1.9296573173537
$ws.Range("BF5").Value = 1.07311885099708
$ws.Range("BG5").Value = 1.80941220814691
$ws.Range("BH5").Value = 1.21888735134514
$ws.Range("BI5").Value = 2.43951142928833
$ws.Range("BJ5").Value = 0.628863222436781
$ws.Range("BK5").Value = 0.951082338779688
$ws.Range("BL5").Value = 0.292596876818483
$ws.Range("BM5").Value = 0.483266420070833
$ws.Range("BN5").Value = 0.962385181778037
$ws.Range("BO5").Value = 0.728212927659969
$ws.Range("BP5").Value = 0.836699456140951
$ws.Range("BQ5").Value = 1.08864605487308
$ws.Range("BR5").Value = 0.446491596043145
$ws.Range("BS5").Value = 0.330390554625628
$ws.Range("BT5").Value = 2.61202665006866
$ws.Range("BU5").Value = 0.275493835123624
$ws.Range("BV5").Value = 1.44045804834025
$ws.Range("BW5").Value = 0.287285179389946
$ws.Range("BX5").Value = 0.43210102735194
$ws.Range("BY5").Value = 0.824702406317234
$ws.Range("BZ5").Value = 0.846591424189426
$ws.Range("CA5").Value = 0.557986251949765
$ws.Range("CB5").Value = 1.23213683284491
$ws.Range("CC5").Value = 0.742444884708347
$ws.Range("CD5").Value = 0.360132811982458
$ws.Range("C6").Value = 0.913430780883414
$ws.Range("D6").Value = 0.659871641927637
$ws.Range("E6").Value = 1.6584351174796
$ws.Range("F6").Value = 0.549286251363516
$ws.Range("G6").Value = 0.329609301642997
$ws.Range("H6").Value = 0.366407171288354
$ws.Range("I6").Value = 0.55908119269249
$ws.Range("J6").Value = 0.622519915515249
$ws.Range("K6").Value = 0.650571355634185
$ws.Range("L6").Value = 1.23588190480388
$ws.Range("M6").Value = 0.639134796594336
$ws.Range("N6").Value = 1.07297375642874
$ws.Range("O6").Value = 1.14163488258397
$ws.Range("P6").Value = 0.152544815012446
$ws.Range("Q6").Value = 0.36620748218017
$ws.Range("R6").Value = 1.69034284922602
$ws.Range("S6").Value = 2.16278726919419
$ws.Range("T6").Value = 0.901385231981898
$ws.Range("U6").Value = 3.89013742108754
$ws.Range("V6").Value = 1.54279116982724
$ws.Range("W6").Value = 1.78124595858254
$ws.Range("X6").Value = 1.07599072788398
$ws.Range("Y6").Value = 0.231635130080127
$ws.Range("Z6").Value = 0.595548299213617
$ws.Range("AA6").Value = 0.476028785007891
$ws.Range("AB6").Value = 1.10923961594865
$ws.Range("AC6").Value = 0.38482119862646
$ws.Range("AD6").Value = 0.104530137666099
$ws.Range("AE6").Value = 0.701288856535181
$ws.Range("AF6").Value = 1.06175830747538
$ws.Range("AG6").Value = 1.1710520526453
$ws.Range("AH6").Value = 0.85264832378249
$ws.Range("AI6").Value = 1.25596961604822
$ws.Range("AJ6").Value = 0.153301294158627
$ws.Range("AK6").Value = 1.7972716240707
$ws.Range("AL6").Value = 2.37695223127668
$ws.Range("AM6").Value = 0.693377933530146
$ws.Range("AN6").Value = 0.519599684945222
$ws.Range("AO6").Value = 0.984891646597838
$ws.Range("AP6").Value = 0.697158993068381
$ws.Range("AQ6").Value = 1.34084080121341
$ws.Range("AR6").Value = 0.45885020088867
$ws.Range("AS6").Value = 1.76828297141319
$ws.Range("AT6").Value = 0.59372616900354
$ws.Range("AU6").Value = 0.248948814493659
$ws.Range("AV6").Value = 0.411621548890464
$ws.Range("AW6").Value = 0.122618092715673
$ws.Range("AX6").Value = 1.23010354784087
$ws.Range("AY6").Value = 0.343453674014832
$ws.Range("AZ6").Value = 0.59782860991356
$ws.Range("BA6").Value = 1.17010359499802
$ws.Range("BB6").Value = 0.239910080141761
$ws.Range("BC6").Value = 1.12373929947633
$ws.Range("BD6").Value = 2.30229025945731
$ws.Range("BE6").Value = 0.0550904218558096
$ws.Range("BF6").Value = 0.829430214624499
$ws.Range("BG6").Value = 0.573043492809938
$ws.Range("BH6").Value = 0.696465044278583
$ws.Range("BI6").Value = 0.888369628439921
$ws.Range("BJ6").Value = 0.204366966224567
$ws.Range("BK6").Value = 0.402294368535011
$ws.Range("BL6").Value = 0.396229857851652
$ws.Range("BM6").Value = 0.250828435964123
$ws.Range("BN6").Value = 0.224682803529888
$ws.Range("BO6").Value = 1.08508802142349
$ws.Range("BP6").Value = 0.772319507694875
$ws.Range("BQ6").Value = 0.219263327709459
$ws.Range("BR6").Value = 1.66756309882958
$ws.Range("BS6").Value = 1.12309878333115
$ws.Range("BT6").Value = 0.425553656299997
$ws.Range("BU6").Value = 0.387701137820165
$ws.Range("BV6").Value = 0.83631257891351
$ws.Range("BW6").Value = 1.57765172204105
$ws.Range("BX6").Value = 0.32105076641665
$ws.Range("BY6").Value = 1.44135267871973
$ws.Range("BZ6").Value = 0.204672409041328
$ws.Range("CA6").Value = 0.916507638940696
$ws.Range("CB6").Value = 1.39568229146966
$ws.Range("CC6").Value = 1.12579505337052
$ws.Range("CD6").Value = 0.229797084625144
$ws.Range("C7").Value = 0.687470634274801
$ws.Range("D7").Value = 1.46766436387023
$ws.Range("E7").Value = 2.18855606011487
$ws.Range("F7").Value = 1.45345582072905
$ws.Range("G7").Value = 1.02791858755953
$ws.Range("H7").Value = 2.52569536055792
$ws.Range("I7").Value = 0.307862986498346
$ws.Range("J7").Value = 1.30397514917735
$ws.Range("K7").Value = 1.09139853437763
$ws.Range("L7").Value = 1.28229535504821
$ws.Range("M7").Value = 0.979158959547126
$ws.Range("N7").Value = 2.45718700018378
$ws.Range("O7").Value = 0.593056693589852
$ws.Range("P7").Value = 1.83927125622366
$ws.Range("Q7").Value = 0.949474095644026
$ws.Range("R7").Value = 1.26627747656772
$ws.Range("S7").Value = 0.628801311647384
$ws.Range("T7").Value = 2.33705024962527
$ws.Range("U7").Value = 1.8511697257789
$ws.Range("V7").Value = 0.489434424738348
$ws.Range("W7").Value = 1.15457324161673
$ws.Range("X7").Value = 1.13153187319078
$ws.Range("Y7").Value = 0.484138931079433
$ws.Range("Z7").Value = 0.480094143840609
$ws.Range("AA7").Value = 0.88531379179075
$ws.Range("AB7").Value = 1.84786856203444
$ws.Range("AC7").Value = 1.5163704134357
$ws.Range("AD7").Value = 1.20487539519469
$ws.Range("AE7").Value = 0.481518388272385
$ws.Range("AF7").Value = 1.61807859052383
$ws.Range("AG7").Value = 0.954893910992704
$ws.Range("AH7").Value = 0.12186295278884
$ws.Range("AI7").Value = 0.347180449186504
$ws.Range("AJ7").Value = 1.35388411656062
$ws.Range("AK7").Value = 2.34911599667082
$ws.Range("AL7").Value = 0.5461364882992
$ws.Range("AM7").Value = 1.10363178139949
$ws.Range("AN7").Value = 0.56688096187637
$ws.Range("AO7").Value = 1.92164536273046
$ws.Range("AP7").Value = 0.813613182121406
$ws.Range("AQ7").Value = 0.598573265178622
$ws.Range("AR7").Value = 0.977468195164536
$ws.Range("AS7").Value = 1.65421488423566
$ws.Range("AT7").Value = 0.234621264834266
$ws.Range("AU7").Value = 1.29880454486761
$ws.Range("AV7").Value = 0.216486269713072
$ws.Range("AW7").Value = 0.392925780152995
$ws.Range("AX7").Value = 1.02499645939834
$ws.Range("AY7").Value = 0.605737347887181
$ws.Range("AZ7").Value = 2.23910803612891
$ws.Range("BA7").Value = 0.978209655534675
$ws.Range("BB7").Value = 0.52288010065481
$ws.Range("BC7").Value = 1.10972537298854
$ws.Range("BD7").Value = 0.0524165742561746
$ws.Range("BE7").Value = 0.315243780974024
$ws.Range("BF7").Value = 1.16998735447494
$ws.Range("BG7").Value = 0.498746391072674
$ws.Range("BH7").Value = 0.956784346254679
$ws.Range("BI7").Value = 0.747268715723711
$ws.Range("BJ7").Value = 0.962807198756881
$ws.Range("BK7").Value = 0.701367020443996
$ws.Range("BL7").Value = 0.580596086038379
$ws.Range("BM7").Value = 1.08867011512743
$ws.Range("BN7").Value = 0.949829743448364
$ws.Range("BO7").Value = 0.25687828595836
$ws.Range("BP7").Value = 1.83358416128635
$ws.Range("BQ7").Value = 2.31259642266001
$ws.Range("BR7").Value = 0.356944748818574
$ws.Range("BS7").Value = 1.0685311228694
$ws.Range("BT7").Value = 1.69607416751229
$ws.Range("BU7").Value = 0.235591746719439
$ws.Range("BV7").Value = 0.582664758875908
$ws.Range("BW7").Value = 1.46195308575025
$ws.Range("BX7").Value = 1.05623219286383
$ws.Range("BY7").Value = 0.422989019600374
$ws.Range("BZ7").Value = 0.269739849174773
$ws.Range("CA7").Value = 1.66268624554193
$ws.Range("CB7").Value = 0.227490641844379
$ws.Range("CC7").Value = 0.744177203655055
$ws.Range("CD7").Value = 0.873712176117503
$ws.Range("C8").Value = 0.779993311680042
$ws.Range("D8").Value = 1.22713679256499
$ws.Range("E8").Value = 1.02140338745146
$ws.Range("F8").Value = 1.19969043805205
$ws.Range("G8").Value = 0.196100933509036
$ws.Range("H8").Value = 0.642106280270053
$ws.Range("I8").Value = 0.833603511417836
$ws.Range("J8").Value = 0.724067093164775
$ws.Range("K8").Value = 2.81152845168912
$ws.Range("L8").Value = 0.319423388438489
$ws.Range("M8").Value = 0.377803332025251
$ws.Range("N8").Value = 0.853555868995608
$ws.Range("O8").Value = 0.373668526893412
$ws.Range("P8").Value = 1.7965386999557
$ws.Range("Q8").Value = 0.0295066231252608
$ws.Range("R8").Value = 1.3856165944888
$ws.Range("S8").Value = 0.923089855705091
$ws.Range("T8").Value = 0.615311226493091
$ws.Range("U8").Value = 0.60846218815891
$ws.Range("V8").Value = 1.71008821251539
$ws.Range("W8").Value = 0.549219875134376
$ws.Range("X8").Value = 0.26901194867354
$ws.Range("Y8").Value = 0.474021041573438
$ws.Range("Z8").Value = 1.69635204499147
$ws.Range("AA8").Value = 1.64260917032037
$ws.Range("AB8").Value = 1.38822883590653
$ws.Range("AC8").Value = 1.70414793518091
$ws.Range("AD8").Value = 1.22830965167684
$ws.Range("AE8").Value = 1.04289971921258
$ws.Range("AF8").Value = 0.54296819845926
$ws.Range("AG8").Value = 0.783317451285665
$ws.Range("AH8").Value = 0.654249718346481
$ws.Range("AI8").Value = 0.117538491280346
$ws.Range("AJ8").Value = 2.00489301285186
$ws.Range("AK8").Value = 0.510714569459694
$ws.Range("AL8").Value = 0.315144881643615
$ws.Range("AM8").Value = 0.298231508407686
$ws.Range("AN8").Value = 0.675793200028313
$ws.Range("AO8").Value = 1.18899654381185
$ws.Range("AP8").Value = 0.893178000673401
$ws.Range("AQ8").Value = 0.494360681078475
$ws.Range("AR8").Value = 1.90628846739295
$ws.Range("AS8").Value = 1.08564937174791
$ws.Range("AT8").Value = 0.316324002916198
$ws.Range("AU8").Value = 1.16309024480222
$ws.Range("AV8").Value = 0.275887148813052
$ws.Range("AW8").Value = 1.53260095038704
$ws.Range("AX8").Value = 0.346177932908771
$ws.Range("AY8").Value = 0.231879306776009
$ws.Range("AZ8").Value = 0.141291092318734
$ws.Range("BA8").Value = 0.558744629874472
$ws.Range("BB8").Value = 0.3869841887994
$ws.Range("BC8").Value = 1.16417528677313
$ws.Range("BD8").Value = 0.463360592938482
$ws.Range("BE8").Value = 1.40416571054711
$ws.Range("BF8").Value = 0.28222775917719
$ws.Range("BG8").Value = 2.37218005021667
$ws.Range("BH8").Value = 0.400802575858037
$ws.Range("BI8").Value = 1.12905521485988
$ws.Range("BJ8").Value = 0.441568976226809
$ws.Range("BK8").Value = 1.85168680070805
$ws.Range("BL8").Value = 0.611783753245756
$ws.Range("BM8").Value = 2.07026739282324
$ws.Range("BN8").Value = 0.589127456764784
$ws.Range("BO8").Value = 1.0737488787304
$ws.Range("BP8").Value = 0.162360747107954
$ws.Range("BQ8").Value = 0.620034112967267
$ws.Range("BR8").Value = 1.25140538937482
$ws.Range("BS8").Value = 0.214442722601544
$ws.Range("BT8").Value = 0.367595334132358
$ws.Range("BU8").Value = 1.14849198608715
$ws.Range("BV8").Value = 2.69835329595507
$ws.Range("BW8").Value = 1.52386036020749
$ws.Range("BX8").Value = 1.14807304183863
$ws.Range("BY8").Value = 0.380291130366733
$ws.Range("BZ8").Value = 1.00727750811354
$ws.Range("CA8").Value = 1.12429687817588
$ws.Range("CB8").Value = 0.795055570083782
$ws.Range("CC8").Value = 0.653854747879739
$ws.Range("CD8").Value = 1.25168576188696
$ws.Range("C9").Value = 0.607725969598161
$ws.Range("D9").Value = 0.256967816864164
$ws.Range("E9").Value = 1.18840380373923
$ws.Range("F9").Value = 1.03573459063886
$ws.Range("G9").Value = 0.354894420515794
$ws.Range("H9").Value = 1.25002133126889
$ws.Range("I9").Value = 0.527915961095092
$ws.Range("J9").Value = 1.86892077859615
$ws.Range("K9").Value = 0.349697420318461
$ws.Range("L9").Value = 0.127920813220847
$ws.Range("M9").Value = 1.90820462626749
$ws.Range("N9").Value = 0.704363067093298
$ws.Range("O9").Value = 0.889286152642572
$ws.Range("P9").Value = 0.922737874343848
$ws.Range("Q9").Value = 0.768871372365159
$ws.Range("R9").Value = 0.346638859604964
$ws.Range("S9").Value = 0.230133238893128
$ws.Range("T9").Value = 1.55612166440216
$ws.Range("U9").Value = 0.162329506729288
$ws.Range("V9").Value = 0.866011907929836
$ws.Range("W9").Value = 1.46685570762755
$ws.Range("X9").Value = 1.60591357727841
$ws.Range("Y9").Value = 0.436094076447862
$ws.Range("Z9").Value = 1.88998734577061
$ws.Range("AA9").Value = 0.115419791689101
$ws.Range("AB9").Value = 3.27990202183751
$ws.Range("AC9").Value = 0.315989476751749
$ws.Range("AD9").Value = 0.651777117783395
$ws.Range("AE9").Value = 0.919356236993446
$ws.Range("AF9").Value = 0.356255950369711
$ws.Range("AG9").Value = 1.5555484798034
$ws.Range("AH9").Value = 0.707458312157688
$ws.Range("AI9").Value = 0.991727114014596
$ws.Range("AJ9").Value = 0.808349967382296
$ws.Range("AK9").Value = 0.815273143796156
$ws.Range("AL9").Value = 1.1824317937347
$ws.Range("AM9").Value = 0.222523906543121
$ws.Range("AN9").Value = 1.61282942319752
$ws.Range("AO9").Value = 0.398909216806682
$ws.Range("AP9").Value = 1.47944611261553
$ws.Range("AQ9").Value = 0.0770741130361868
$ws.Range("AR9").Value = 0.539528615726297
$ws.Range("AS9").Value = 1.64795565499245
$ws.Range("AT9").Value = 1.91882020953806
$ws.Range("AU9").Value = 1.03529705339596
$ws.Range("AV9").Value = 0.607389464885687
$ws.Range("AW9").Value = 1.23239388141012
$ws.Range("AX9").Value = 0.458772251054794
$ws.Range("AY9").Value = 1.40250938609649
$ws.Range("AZ9").Value = 1.60891426460527
$ws.Range("BA9").Value = 1.28433623473007
$ws.Range("BB9").Value = 0.169747501740303
$ws.Range("BC9").Value = 1.39031233400662
$ws.Range("BD9").Value = 1.49107914191432
$ws.Range("BE9").Value = 0.794272444676741
$ws.Range("BF9").Value = 1.53658355081364
$ws.Range("BG9").Value = 0.853735349831292
$ws.Range("BH9").Value = 2.71167960702326
$ws.Range("BI9").Value = 1.30119151145236
$ws.Range("BJ9").Value = 1.27000143332568
$ws.Range("BK9").Value = 0.619550438773352
$ws.Range("BL9").Value = 1.51826287421087
$ws.Range("BM9").Value = 1.33165338716946
$ws.Range("BN9").Value = 1.52595566602809
$ws.Range("BO9").Value = 1.44237167816717
$ws.Range("BP9").Value = 0.907976983239359
$ws.Range("BQ9").Value = 1.1063657483934
$ws.Range("BR9").Value = 0.907955541450917
$ws.Range("BS9").Value = 0.849403470166042
$ws.Range("BT9").Value = 1.55387928199571
$ws.Range("BU9").Value = 0.528014267168349
$ws.Range("BV9").Value = 0.284338985854892
$ws.Range("BW9").Value = 0.283721700570164
$ws.Range("BX9").Value = 2.80695775899244
$ws.Range("BY9").Value = 0.243855621268796
$ws.Range("BZ9").Value = 1.10139032256875
$ws.Range("CA9").Value = 1.12657678218989
$ws.Range("CB9").Value = 1.20052370243686
$ws.Range("CC9").Value = 1.05103462945012
$ws.Range("CD9").Value = 1.17029345003651
$ws.Range("C10").Value = 0.316928726009427
$ws.Range("D10").Value = 0.0728307356764099
$ws.Range("E10").Value = 2.47054439694463
$ws.Range("F10").Value = 1.32316508917093
$ws.Range("G10").Value = 0.935995921255251
$ws.Range("H10").Value = 0.995690194719253
$ws.Range("I10").Value = 1.41095626919158
$ws.Range("J10").Value = 1.96960203466538
$ws.Range("K10").Value = 0.642617949471252
$ws.Range("L10").Value = 0.0294855281616106
$ws.Range("M10").Value = 0.488736803638604
$ws.Range("N10").Value = 0.964261774808223
$ws.Range("O10").Value = 0.142685742076882
$ws.Range("P10").Value = 0.813707466875029
$ws.Range("Q10").Value = 0.872000510793204
$ws.Range("R10").Value = 0.563766657916921
$ws.Range("S10").Value = 0.371835890836891
$ws.Range("T10").Value = 0.483469117287807
$ws.Range("U10").Value = 0.554256331751049
$ws.Range("V10").Value = 0.120525401881687
$ws.Range("W10").Value = 0.767491443355059
$ws.Range("X10").Value = 0.837267223687003
$ws.Range("Y10").Value = 0.246188772088726
$ws.Range("Z10").Value = 0.396879888280418
$ws.Range("AA10").Value = 0.295010032543904
$ws.Range("AB10").Value = 1.17982845538142
$ws.Range("AC10").Value = 2.99630375503952
$ws.Range("AD10").Value = 0.582629269617024
$ws.Range("AE10").Value = 1.00609111627423
$ws.Range("AF10").Value = 0.183217490927327
$ws.Range("AG10").Value = 1.02926430010419
$ws.Range("AH10").Value = 0.940511907367985
$ws.Range("AI10").Value = 0.103290203933459
$ws.Range("AJ10").Value = 0.647455030808851
$ws.Range("AK10").Value = 1.64670493329013
$ws.Range("AL10").Value = 0.296072826067226
$ws.Range("AM10").Value = 0.190637662779299
$ws.Range("AN10").Value = 0.559631573191953
$ws.Range("AO10").Value = 1.61291467415221
$ws.Range("AP10").Value = 1.94573874931194
$ws.Range("AQ10").Value = 1.3250811702316
$ws.Range("AR10").Value = 1.17986876365372
$ws.Range("AS10").Value = 1.04466067500861
$ws.Range("AT10").Value = 2.53842921142866
$ws.Range("AU10").Value = 0.606628965363358
$ws.Range("AV10").Value = 0.397604143708638
$ws.Range("AW10").Value = 0.432842448825171
$ws.Range("AX10").Value = 0.31504809388954
$ws.Range("AY10").Value = 1.19797663616955
$ws.Range("AZ10").Value = 0.77726536327133
$ws.Range("BA10").Value = 0.997311211536098
$ws.Range("BB10").Value = 0.410291987543381
$ws.Range("BC10").Value = 0.768205872783498
$ws.Range("BD10").Value = 0.462930027702256
$ws.Range("BE10").Value = 1.0734777675058
$ws.Range("BF10").Value = 0.740900170341566
$ws.Range("BG10").Value = 0.392970119511369
$ws.Range("BH10").Value = 0.273445672340824
$ws.Range("BI10").Value = 1.1822609010457
$ws.Range("BJ10").Value = 0.852382508989545
$ws.Range("BK10").Value = 0.468996567101669
$ws.Range("BL10").Value = 1.55450461723018
$ws.Range("BM10").Value = 3.61523234682548
$ws.Range("BN10").Value = 0.273731677073769
$ws.Range("BO10").Value = 0.420510822702785
$ws.Range("BP10").Value = 0.263086415238878
$ws.Range("BQ10").Value = 1.11167962068122
$ws.Range("BR10").Value = 0.884401846364168
$ws.Range("BS10").Value = 0.989209734795689
$ws.Range("BT10").Value = 1.5720914325048
$ws.Range("BU10").Value = 0.609223329215422
$ws.Range("BV10").Value = 0.669489910153273
$ws.Range("BW10").Value = 0.645484343085754
$ws.Range("BX10").Value = 1.28419720554329
$ws.Range("BY10").Value = 1.8138865441747
$ws.Range("BZ10").Value = 3.99428733170229
$ws.Range("CA10").Value = 0.564867728097179
$ws.Range("CB10").Value = 1.61620434590235
$ws.Range("CC10").Value = 0.637809659381476
$ws.Range("CD10").Value = 0.716792904304481
$ws.Range("C11").Value = 0.220978098239653
$ws.Range("D11").Value = 0.207769101296031
$ws.Range("E11").Value = 1.73794428688282
$ws.Range("F11").Value = 1.61923871438846
$ws.Range("G11").Value = 0.311992508254432
$ws.Range("H11").Value = 1.36299267133993
$ws.Range("I11").Value = 0.396153067834607
$ws.Range("J11").Value = 0.588856278000165
$ws.Range("K11").Value = 1.85377825423113
$ws.Range("L11").Value = 0.430357272229962
$ws.Range("M11").Value = 1.00499867836694
$ws.Range("N11").Value = 0.417208165994508
$ws.Range("O11").Value = 1.02452931989108
$ws.Range("P11").Value = 0.711799006869154
$ws.Range("Q11").Value = 0.515435847192503
$ws.Range("R11").Value = 0.54280654508543
$ws.Range("S11").Value = 2.72668860237822
$ws.Range("T11").Value = 1.02198177775651
$ws.Range("U11").Value = 0.173929553880298
$ws.Range("V11").Value = 0.0995274308924302
$ws.Range("W11").Value = 0.867279528091831
$ws.Range("X11").Value = 0.968363294535982
$ws.Range("Y11").Value = 0.279295295368806
$ws.Range("Z11").Value = 1.78223714334966
$ws.Range("AA11").Value = 1.16392399348288
$ws.Range("AB11").Value = 1.21144030338683
$ws.Range("AC11").Value = 1.64542481950428
$ws.Range("AD11").Value = 0.330997044472793
$ws.Range("AE11").Value = 0.722945581192247
$ws.Range("AF11").Value = 0.577528208888047
$ws.Range("AG11").Value = 0.502199637183231
$ws.Range("AH11").Value = 0.557312874523984
$ws.Range("AI11").Value = 0.422247272989754
$ws.Range("AJ11").Value = 0.803024341283368
$ws.Range("AK11").Value = 0.844761239987158
$ws.Range("AL11").Value = 0.0600380303556005
$ws.Range("AM11").Value = 1.08250815967878
$ws.Range("AN11").Value = 0.447405633467945
$ws.Range("AO11").Value = 0.500554885595155
$ws.Range("AP11").Value = 1.04574120780773
$ws.Range("AQ11").Value = 1.28373127342109
$ws.Range("AR11").Value = 1.03449689195515
$ws.Range("AS11").Value = 2.40244916255559
$ws.Range("AT11").Value = 0.265917731056268
$ws.Range("AU11").Value = 2.39516621106476
$ws.Range("AV11").Value = 1.52748341893591
$ws.Range("AW11").Value = 0.746396492473633
$ws.Range("AX11").Value = 0.395084721114261
$ws.Range("AY11").Value = 0.29947543415625
$ws.Range("AZ11").Value = 1.12705363437738
$ws.Range("BA11").Value = 0.151988430412107
$ws.Range("BB11").Value = 0.511693098869842
$ws.Range("BC11").Value = 0.202710958231061
$ws.Range("BD11").Value = 0.0915565812034053
$ws.Range("BE11").Value = 1.52809606371176
$ws.Range("BF11").Value = 0.937408277035619
$ws.Range("BG11").Value = 0.360011060203765
$ws.Range("BH11").Value = 0.934722849549275
$ws.Range("BI11").Value = 2.42116982837018
$ws.Range("BJ11").Value = 0.575870582884403
$ws.Range("BK11").Value = 1.72796226469874
$ws.Range("BL11").Value = 0.868319987211924
$ws.Range("BM11").Value = 0.481593573236081
$ws.Range("BN11").Value = 3.21643848490712
$ws.Range("BO11").Value = 1.57724583004414
$ws.Range("BP11").Value = 0.0560659618207038
$ws.Range("BQ11").Value = 1.13787529232522
$ws.Range("BR11").Value = 0.655583330172966
$ws.Range("BS11").Value = 0.375026722595728
$ws.Range("BT11").Value = 0.639449634479694
$ws.Range("BU11").Value = 1.32003597894287
$ws.Range("BV11").Value = 1.60825514529555
$ws.Range("BW11").Value = 1.54777144423028
$ws.Range("BX11").Value = 0.283530812583339
$ws.Range("BY11").Value = 3.21292555732049
$ws.Range("BZ11").Value = 0.624633069406631
$ws.Range("CA11").Value = 0.278144099341691
$ws.Range("CB11").Value = 0.832676594968378
$ws.Range("CC11").Value = 1.28021145059423
$ws.Range("CD11").Value = 3.90685060840999
